$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for 73c57d2d... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-07 11:03:42"

# "zh-cn" sheet: Correspond Handoff/Handback Datetime for 73c57d2d... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-07 11:03:37"
$wsZhCn.Range("K4").Value = "2016-09-07 11:03:58"

# "de-de" sheet: Correspond Handback Datetime for 73c57d2d... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-09-07 11:04:16"
